$d = $word.ActiveDocument

$paragraphs = @($d.Paragraphs)
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    $p = $paragraphs[$i]
    $t = $p.Range.Text.Trim()
    if ($t -eq "Bild vertikales Gleis reparieren") {
        # Build a range spanning this paragraph and the following (empty) one,
        # then delete them both in one go.
        $start = $p.Range.Start
        $next = $paragraphs[$i + 1]
        $end = $next.Range.End
        $r = $d.Range($start, $end)
        $r.Delete()
        break
    }
}
